$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 241, shifting rows 241:297 down to 243:299
$ws.Range("A241:R242").EntireRow.Insert()

# Fill the new row 241 (Primera)
$ws.Cells.Item(241, 1).Value = 9
$ws.Cells.Item(241, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(241, 3).Value = "Metropolitana"
$ws.Cells.Item(241, 4).Value = 44736
$ws.Cells.Item(241, 5).Value = 13
$ws.Cells.Item(241, 6).Value = 100112017
$ws.Cells.Item(241, 7).Value = "Apio"
$ws.Cells.Item(241, 8).Value = "Americana (o)"
$ws.Cells.Item(241, 9).Value = "Primera"
$ws.Cells.Item(241, 10).Value = 52
$ws.Cells.Item(241, 11).Value = 6000
$ws.Cells.Item(241, 12).Value = 7000
$ws.Cells.Item(241, 13).Value = 6500
$ws.Cells.Item(241, 14).Value = "`$/docena de matas"
$ws.Cells.Item(241, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(241, 16).Value = 1083
$ws.Cells.Item(241, 17).Value = 6
$ws.Cells.Item(241, 18).Value = "Hortaliza"

# Fill the new row 242 (Segunda)
$ws.Cells.Item(242, 1).Value = 9
$ws.Cells.Item(242, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(242, 3).Value = "Metropolitana"
$ws.Cells.Item(242, 4).Value = 44736
$ws.Cells.Item(242, 5).Value = 13
$ws.Cells.Item(242, 6).Value = 100112017
$ws.Cells.Item(242, 7).Value = "Apio"
$ws.Cells.Item(242, 8).Value = "Americana (o)"
$ws.Cells.Item(242, 9).Value = "Segunda"
$ws.Cells.Item(242, 10).Value = 34
$ws.Cells.Item(242, 11).Value = 5000
$ws.Cells.Item(242, 12).Value = 5000
$ws.Cells.Item(242, 13).Value = 5000
$ws.Cells.Item(242, 14).Value = "`$/docena de matas"
$ws.Cells.Item(242, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(242, 16).Value = 833
$ws.Cells.Item(242, 17).Value = 6
$ws.Cells.Item(242, 18).Value = "Hortaliza"
